$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N ("Late") on the Repayment
# Schedule sheet, shifting the existing N/O/P columns ("Late", blank,
# "Outstanding") one place to the right (-> O/P/Q).
$ws.Columns("N").Insert() | Out-Null

# Make "Repayment Schedule" the active sheet/tab (was "Transactions")
# and leave the selection on R7.
$ws.Select() | Out-Null
$ws.Range("R7").Select() | Out-Null
